$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'26.980.28"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.47%  '
$ws.Range("D3").Value = "'1.826.16"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.15%  '
$ws.Range("E4").Value = '  -0.29%  '
$ws.Range("D5").Value = "'312.10"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.16%  '
$ws.Range("D7").Value = "'0.4626"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  +1.75%  '
$ws.Range("D9").Value = "'0.07331"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.55%  '
$ws.Range("D10").Value = "'0.8748"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.51%  '
$ws.Range("D11").Value = "'0.07940"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +4.40%  '
$ws.Range("E12").Value = '  -1.76%  '
$ws.Range("D13").Value = "'1.893.17"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +1.59%  '
$ws.Range("D14").Value = "'5.337"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.25%  '
$ws.Range("D15").Value = "'6.534"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.88%  '
$ws.Range("D16").Value = "'91.24"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.48%  '
$ws.Range("D17").Value = "'1.007"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.16%  '
$ws.Range("D18").Value = "'0.000008872"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +2.64%  '
$ws.Range("D19").Value = "'1.005"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.36%  '
$ws.Range("D20").Value = "'14.78"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("D21").Value = "'26.920.46"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.67%  '
$ws.Range("D22").Value = "'5.102"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.89%  '
$ws.Range("D23").Value = "'10.53"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.30%  '
$ws.Range("D24").Value = "'2.137.70"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.33%  '
$ws.Range("D25").Value = "'153.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.04%  '
$ws.Range("D26").Value = "'1.849"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.26%  '
$ws.Range("D27").Value = "'18.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.44%  '
$ws.Range("D28").Value = "'2.036"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.16%  '
$ws.Range("D29").Value = "'5.140"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.27%  '
$ws.Range("D30").Value = "'115.59"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.58%  '
$ws.Range("D31").Value = "'0.08903"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.19%  '
$ws.Range("D33").Value = "'0.7276"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.99%  '
$ws.Range("D34").Value = "'4.435"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.34%  '
$ws.Range("D35").Value = "'1.130"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.92%  '
$ws.Range("D36").Value = "'2.491"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.78%  '
$ws.Range("D37").Value = "'0.01954"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.67%  '
$ws.Range("D38").Value = "'1.068"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.56%  '
$ws.Range("D39").Value = "'0.05224"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.73%  '
$ws.Range("D40").Value = "'2.942"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.27%  '
$ws.Range("D41").Value = "'7.090"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.07%  '
$ws.Range("D42").Value = "'0.5156"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.18%  '
$ws.Range("D43").Value = "'0.1623"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.61%  '
$ws.Range("B44").Value = 'Decentraland'
$ws.Range("C44").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D44").Value = "'0.4847"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.04%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = "'8.175"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.16%  '
$ws.Range("E46").Value = '  -0.29%  '
$ws.Range("D47").Value = "'10.17"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.27%  '
$ws.Range("D48").Value = "'102.74"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.31%  '
$ws.Range("E49").Value = '  -0.15%  '
$ws.Range("D50").Value = "'0.06187"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.07%  '
$ws.Range("D51").Value = "'64.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.14%  '
